# Weekly update: insert one new price record before the existing row 70,
# pushing the old rows 70-196 down to 71-197 (dimension grows from
# A1:R196 to A1:R197), then populate the newly-inserted row with the
# new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70 (existing rows 70.. shift down by 1)
$ws.Rows.Item(70).Insert()

# Fill in the new row 70 with this week's record
$ws.Cells.Item(70, 1).Value  = 10
$ws.Cells.Item(70, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value  = "La Araucanía"
$ws.Cells.Item(70, 4).Value  = 44469
$ws.Cells.Item(70, 5).Value  = 9
$ws.Cells.Item(70, 6).Value  = 100112009
$ws.Cells.Item(70, 7).Value  = "Acelga"
$ws.Cells.Item(70, 8).Value  = "Sin especificar"
$ws.Cells.Item(70, 9).Value  = "Primera"
$ws.Cells.Item(70, 10).Value = 110
$ws.Cells.Item(70, 11).Value = 7000
$ws.Cells.Item(70, 12).Value = 8000
$ws.Cells.Item(70, 13).Value = 7455
$ws.Cells.Item(70, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(70, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(70, 16).Value = 621
$ws.Cells.Item(70, 17).Value = 12
$ws.Cells.Item(70, 18).Value = "Hortaliza"
